# Update GUI. Added radio buttons.
#
# - Typography sheet: add a new font row ("RB_Indicators", arial.ttf,
#   size 11, bpp 4) in the first free row of the Typography table.
# - Translation sheet: update an existing text's base (GB) value, and add
#   the new Text ID rows that drive the new radio-button indicator texts
#   (Run/Stop, Setpoint, Program) plus two more generic value rows.

$wb = $excel.ActiveWorkbook
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------
# Typography sheet - new Font row (row 10)
# ---------------------------------------------------------------------
$wsTypography.Range("B10").Value = "RB_Indicators"
$wsTypography.Range("C10").Value = "arial.ttf"
$wsTypography.Range("D10").Value = 11
$wsTypography.Range("E10").Value = 4
$wsTypography.Range("F10").Value = "?"

# Keep these cells on the sheet's normal/default style (matches the
# other data rows in the table, which carry no explicit cell style).
$wsTypography.Range("B10:J10").Style = "Normal"

# ---------------------------------------------------------------------
# Translation sheet - new Text ID rows (rows 41-48)
# ---------------------------------------------------------------------
$wsTranslation.Range("B41").Value = "SingleUseId112"
$wsTranslation.Range("C41").Value = "t_AI"
$wsTranslation.Range("D41").Value = "Left"
$wsTranslation.Range("E41").Value = "<value>"
$wsTranslation.Range("F41").Value = "LTR"

$wsTranslation.Range("B42").Value = "SingleUseId113"
$wsTranslation.Range("C42").Value = "t_AI"
$wsTranslation.Range("D42").Value = "Left"
# Force text storage - "0" would otherwise be auto-coerced to a number.
$wsTranslation.Range("E42").NumberFormat = "@"
$wsTranslation.Range("E42").Value = "0"
$wsTranslation.Range("F42").Value = "LTR"

$wsTranslation.Range("B43").Value = "SingleUseId114"
$wsTranslation.Range("C43").Value = "RB_Indicators"
$wsTranslation.Range("D43").Value = "Center"
$wsTranslation.Range("E43").Value = "<value> "
$wsTranslation.Range("F43").Value = "LTR"

$wsTranslation.Range("B44").Value = "SingleUseId115"
$wsTranslation.Range("C44").Value = "RB_Indicators"
$wsTranslation.Range("D44").Value = "Left"
$wsTranslation.Range("E44").Value = "RS"
$wsTranslation.Range("F44").Value = "LTR"

$wsTranslation.Range("B45").Value = "SingleUseId116"
$wsTranslation.Range("C45").Value = "RB_Indicators"
$wsTranslation.Range("D45").Value = "Center"
$wsTranslation.Range("E45").Value = "<value> "
$wsTranslation.Range("F45").Value = "LTR"

$wsTranslation.Range("B46").Value = "SingleUseId117"
$wsTranslation.Range("C46").Value = "RB_Indicators"
$wsTranslation.Range("D46").Value = "Left"
$wsTranslation.Range("E46").Value = "SD"
$wsTranslation.Range("F46").Value = "LTR"

$wsTranslation.Range("B47").Value = "SingleUseId118"
$wsTranslation.Range("C47").Value = "RB_Indicators"
$wsTranslation.Range("D47").Value = "Center"
$wsTranslation.Range("E47").Value = "<value> "
$wsTranslation.Range("F47").Value = "LTR"

$wsTranslation.Range("B48").Value = "SingleUseId119"
$wsTranslation.Range("C48").Value = "RB_Indicators"
$wsTranslation.Range("D48").Value = "Left"
$wsTranslation.Range("E48").Value = "PRG"
$wsTranslation.Range("F48").Value = "LTR"

# Match the default (un-styled) look of the rest of the Text ID table.
$wsTranslation.Range("B41:F48").Style = "Normal"

# ---------------------------------------------------------------------
# Translation sheet - existing row 5 ("t_AI" base value) text tweak
# ---------------------------------------------------------------------
$wsTranslation.Range("E5").Value = " :  <value>"
$wsTranslation.Range("E5").Value = " : <value>"
